# Update "想去人数" (interested-count) figures across all worksheets.
# Mirrors the data refresh captured in the commit "Update gh-pages to
# output generated at 456a3b4" (the site's scraper numbers ticked up
# slightly for several events).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1897
$ws1.Range("F3").Value  = 1522
$ws1.Range("F4").Value  = 889
$ws1.Range("F5").Value  = 785
$ws1.Range("F6").Value  = 13386
$ws1.Range("F7").Value  = 13248
$ws1.Range("F11").Value = 567
$ws1.Range("F13").Value = 687
$ws1.Range("F14").Value = 2102
$ws1.Range("F17").Value = 77
$ws1.Range("F19").Value = 402
$ws1.Range("F20").Value = 275
$ws1.Range("F21").Value = 293
$ws1.Range("F22").Value = 428

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 127
$ws2.Range("F9").Value = 33

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 52

# Sheet "全部类型" (All types - combined listing)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1897
$ws4.Range("F4").Value  = 1522
$ws4.Range("F5").Value  = 889
$ws4.Range("F7").Value  = 785
$ws4.Range("F8").Value  = 13386
$ws4.Range("F9").Value  = 13248
$ws4.Range("F13").Value = 567
$ws4.Range("F15").Value = 687
$ws4.Range("F18").Value = 2102
$ws4.Range("F21").Value = 77
$ws4.Range("F25").Value = 52
$ws4.Range("F26").Value = 402
$ws4.Range("F27").Value = 275
$ws4.Range("F28").Value = 293
$ws4.Range("F29").Value = 428
$ws4.Range("F31").Value = 127
$ws4.Range("F34").Value = 33
